$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new snippet rows describing the new line-shape APIs
$ws.Range("A168").Value = "ShapeCollection"
$ws.Range("B168").Value = "addLine"
$ws.Range("C168").Value = "excel-shape-lines"
$ws.Range("D168").Value = "addStraightLine"

$ws.Range("A169").Value = "Shape"
$ws.Range("B169").Value = "line"
$ws.Range("C169").Value = "excel-shape-lines"
$ws.Range("D169").Value = "arrowLine"

# Grow the "Snippets" table to include the newly added rows
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:D169"))

# Update the active selection to reflect where the user ended up after the edit
$ws.Range("A170").Select()
